# Generate Report for Archive
#
# The three entries "96e8afa5-...", "cf667357-..." and "d2660411-..." (rows 5-7 on
# every sheet) get re-sorted: "d2660411-..." now reports status "In Translation"
# and moves to the top of the group (row 5), while "96e8afa5-..." and
# "cf667357-..." shift down one row each but otherwise keep their own data.
# Hyperlink display text follows the cell it is attached to.

$wb = $excel.ActiveWorkbook

function Get-HyperlinkSnapshot {
    param($ws)
    # Snapshot display-text -> Hyperlink object BEFORE any edits so later lookups
    # are unambiguous even once several cells end up sharing text momentarily.
    $map = @{}
    foreach ($h in $ws.Hyperlinks) {
        $map[$h.TextToDisplay] = $h
    }
    return $map
}

function Set-CellText {
    param($ws, [string]$cellRef, [string]$newValue, $hlSnapshot)
    $range = $ws.Range($cellRef)
    $oldValue = $range.Value2
    $range.Value = $newValue
    if ($hlSnapshot -ne $null -and $hlSnapshot.ContainsKey($oldValue)) {
        $hlSnapshot[$oldValue].TextToDisplay = $newValue
    }
}

# ---------- Sheet "Overview" ----------
$ws = $wb.Worksheets.Item("Overview")
$hl = Get-HyperlinkSnapshot $ws

Set-CellText $ws "A5" "d2660411-d592-436e-84cb-7090cc94fda8.md" $hl
Set-CellText $ws "B5" "In Translation" $hl
Set-CellText $ws "C5" "In Translation" $hl
Set-CellText $ws "D5" "2016-03-22 00:35:57" $hl

Set-CellText $ws "A6" "96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.md" $hl
Set-CellText $ws "B6" "Ready for handoff" $hl
Set-CellText $ws "C6" "Ready for handoff" $hl
Set-CellText $ws "D6" "2016-03-22 00:33:37" $hl

Set-CellText $ws "A7" "cf667357-71fc-4ed0-8b1d-8d1dc74b41dd.md" $hl
Set-CellText $ws "B7" "Ready for handoff" $hl
Set-CellText $ws "C7" "Ready for handoff" $hl
Set-CellText $ws "D7" "2016-03-22 00:36:31" $hl

# ---------- Sheet "zh-cn" ----------
$ws = $wb.Worksheets.Item("zh-cn")
$hl = Get-HyperlinkSnapshot $ws

Set-CellText $ws "A5" "d2660411-d592-436e-84cb-7090cc94fda8.md" $hl
Set-CellText $ws "C5" "In Translation" $hl
Set-CellText $ws "D5" "d2660411-d592-436e-84cb-7090cc94fda8.7cecc8beac20682be3c31762b1fed381e7ddd62b.zh-cn.xlf" $hl
Set-CellText $ws "E5" "2016-03-22 00:35:53" $hl

Set-CellText $ws "A6" "96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.md" $hl
Set-CellText $ws "C6" "Ready for handoff" $hl
Set-CellText $ws "D6" "96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.0012e40d796e5c6f54b3c87d5af7bf616b8ae37b.zh-cn.xlf" $hl
Set-CellText $ws "E6" "2016-03-22 00:33:33" $hl

Set-CellText $ws "A7" "cf667357-71fc-4ed0-8b1d-8d1dc74b41dd.md" $hl
Set-CellText $ws "C7" "Ready for handoff" $hl
Set-CellText $ws "D7" "cf667357-71fc-4ed0-8b1d-8d1dc74b41dd.bb122b9ccdade679a9783ff69492a289cd8dd1fb.zh-cn.xlf" $hl
Set-CellText $ws "E7" "2016-03-22 00:36:27" $hl

# ---------- Sheet "de-de" ----------
$ws = $wb.Worksheets.Item("de-de")
$hl = Get-HyperlinkSnapshot $ws

Set-CellText $ws "A5" "d2660411-d592-436e-84cb-7090cc94fda8.md" $hl
Set-CellText $ws "C5" "In Translation" $hl
Set-CellText $ws "D5" "d2660411-d592-436e-84cb-7090cc94fda8.7cecc8beac20682be3c31762b1fed381e7ddd62b.de-de.xlf" $hl
Set-CellText $ws "E5" "2016-03-22 00:35:57" $hl

Set-CellText $ws "A6" "96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.md" $hl
Set-CellText $ws "C6" "Ready for handoff" $hl
Set-CellText $ws "D6" "96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.0012e40d796e5c6f54b3c87d5af7bf616b8ae37b.de-de.xlf" $hl
Set-CellText $ws "E6" "2016-03-22 00:33:37" $hl

Set-CellText $ws "A7" "cf667357-71fc-4ed0-8b1d-8d1dc74b41dd.md" $hl
Set-CellText $ws "C7" "Ready for handoff" $hl
Set-CellText $ws "D7" "cf667357-71fc-4ed0-8b1d-8d1dc74b41dd.bb122b9ccdade679a9783ff69492a289cd8dd1fb.de-de.xlf" $hl
Set-CellText $ws "E7" "2016-03-22 00:36:31" $hl

$wb.Save()
